# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.723.50"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.133.38"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  +2.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.412"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +12.61%  "

$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.131.34"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.751"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.16%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("E13").Value = "  +5.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.356.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.54%  "

$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.725.77"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.108.30"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.27%  "

$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000209"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.15%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "452.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.83%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.45"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.302.53"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +9.35%  "

$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.228"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.33"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.21%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.41%  "

$ws.Range("E36").Value = "  -3.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.47"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.56%  "

$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "485.87"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.65%  "

$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.65%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.63%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.439"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.59%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("E47").Value = "  +3.76%  "

$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("E49").Value = "  +4.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.17%  "

